$d = $word.ActiveDocument

function Insert-NewParaBefore($paraIndex, $text) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $r.Collapse(1)  # wdCollapseStart
    $r.InsertParagraphBefore()
    $newRange = $d.Paragraphs.Item($paraIndex).Range
    $newRange.Text = $text
    $newRange.Font.Size = 14
}

# 1) "Good communicator..." before the empty paragraph preceding "Producer: Henry" (orig para 2)
Insert-NewParaBefore 2 "Good communicator, good at getting team to get to work, good leader"

# After the previous insert, everything shifted by +1.
# 2) "Organized, keen to keep notes..." before the empty paragraph preceding "Art Lead: Xander" (orig para 4 -> now 5)
Insert-NewParaBefore 5 "Organized, keen to keep notes, good communicator"

# shifted by +1 again
# 3) "Knowledgeable about sprites" before empty paragraph preceding "Design Lead: Jun" (orig para 6 -> now 8)
Insert-NewParaBefore 8 "Knowledgeable about sprites"

# shifted by +1 again
# 4) "Picky, attentive, specific on game details " before empty paragraph preceding "Sound Lead: Henry" (orig para 8 -> now 11)
Insert-NewParaBefore 11 "Picky, attentive, specific on game details "

# shifted by +1 again
# 5) "Well versed in sound design/audio editing" before empty paragraph preceding "Tech Lead: Koben" (orig para 10 -> now 14)
Insert-NewParaBefore 14 "Well versed in sound design/audio editing"

# Now handle the Tech Lead paragraph: split "Tech Lead: Koben" into two runs with proofErr tags,
# and move the bookmark to a brand new trailing paragraph.
$n = $d.Paragraphs.Count
$techPara = $d.Paragraphs.Item($n)
$full = $techPara.Range
$innerEnd = $full.End - 1   # exclude the paragraph mark
$inner = $d.Range($full.Start, $innerEnd)

$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:rPr><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve">Tech Lead: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="36"/></w:rPr><w:t>Koben</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$inner.InsertXML($xml)

# Remove the old bookmark that now wraps the split runs (it will be re-added on the new trailing paragraph)
$d.Bookmarks.Item("_GoBack").Delete()

# Insert the final new paragraph after the Tech Lead paragraph, with the bookmark.
$techPara2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$rEnd = $techPara2.Range
$rEnd.Collapse(0)  # wdCollapseEnd
$rEnd.InsertParagraphAfter()
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$lastRange.Text = "Good knowledge of unity and how it interacts with C# code, knowledgeable about C# code"
$lastRange.Font.Size = 14

$lastPara2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$endOfLast = $lastPara2.Range
$endOfLast.Collapse(0)
$d.Bookmarks.Add("_GoBack", $endOfLast)
